# "Updated cryptos list ... with GitHub Actions" -- refresh coin rows
# (name/link/price/1h-volume) from the latest scrape. The source feed
# pushes a new row in at #11 (OKB) and shifts every following coin down
# by one, so most rows below it get a new coin's data; the former last
# row (Aave) drops off the bottom of the 50-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D="28.438.96"; E="  +1.12%  "},
    @{Row=3; D="1.880.07"; E="  +0.27%  "},
    @{Row=4; E="  +1.22%  "},
    @{Row=5; D="316.20"; E="  +1.03%  "},
    @{Row=6; D="1.014"; E="  +1.26%  "},
    @{Row=7; D="0.5142"; E="  +0.16%  "},
    @{Row=8; D="0.3945"; E="  +1.35%  "},
    @{Row=9; D="0.08332"; E="  -0.75%  "},
    @{Row=10; D="1.122"; E="  +0.54%  "},
    @{Row=11; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="41.99"},
    @{Row=12; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="6.277"; E="  +1.05%  "},
    @{Row=13; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.865.28"; E="  -0.62%  "},
    @{Row=14; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="20.39"; E="  -1.61%  "},
    @{Row=15; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.261"; E="  -0.39%  "},
    @{Row=16; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.013"; E="  +1.01%  "},
    @{Row=17; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.00001108"; E="  -0.14%  "},
    @{Row=18; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="91.61"; E="  +0.69%  "},
    @{Row=19; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.06734"; E="  +1.25%  "},
    @{Row=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="17.77"; E="  +0.36%  "},
    @{Row=21; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.014"; E="  +1.15%  "},
    @{Row=22; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.989"; E="  -1.11%  "},
    @{Row=23; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="28.468.21"; E="  +1.06%  "},
    @{Row=24; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="11.16"; E="  +0.16%  "},
    @{Row=25; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.254"; E="  +0.15%  "},
    @{Row=26; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.082.90"; E="  +0.04%  "},
    @{Row=27; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="161.47"; E="  +1.95%  "},
    @{Row=28; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="20.81"; E="  +0.93%  "},
    @{Row=29; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.421"; E="  -3.17%  "},
    @{Row=30; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="127.36"; E="  +1.81%  "},
    @{Row=31; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.1060"; E="  -0.27%  "},
    @{Row=32; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.042"; E="  +0.13%  "},
    @{Row=33; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.876"; E="  -0.16%  "},
    @{Row=34; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.642"; E="  +1.19%  "},
    @{Row=35; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02453"; E="  +0.45%  "},
    @{Row=36; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.06530"; E="  -0.25%  "},
    @{Row=37; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="9.152"; E="  -6.17%  "},
    @{Row=38; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.2189"; E="  +0.14%  "},
    @{Row=39; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.260"; E="  +2.69%  "},
    @{Row=40; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.6478"; E="  -0.55%  "},
    @{Row=41; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.189"; E="  -1.74%  "},
    @{Row=42; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="4.999"; E="  -0.53%  "},
    @{Row=43; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="11.14"; E="  -1.65%  "},
    @{Row=44; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.6050"; E="  -0.92%  "},
    @{Row=45; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="13.12"; E="  +0.02%  "},
    @{Row=46; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="3.703"; E="  +0.80%  "},
    @{Row=47; B="WEMIXTOKEN"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="1.284"; E="  +0.26%  "},
    @{Row=48; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="2.012"; E="  +0.01%  "},
    @{Row=49; B="EOS"; C="https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D="1.213"; E="  -0.35%  "},
    @{Row=50; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="121.70"; E="  +0.04%  "},
    @{Row=51; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.06904"; E="  +0.33%  "}
)

# Cells that look like a plain number (e.g. "41.99") would otherwise be
# auto-converted from text to a numeric value by Excel when assigned via
# .Value, which would also silently rewrite things like "316.20" -> 316.2
# or "0.06904" -> 6.904E-02. Prefixing with an apostrophe forces Excel to
# keep (and store) them as text, matching the source data.
$numericPattern = '^[+-]?(\d+(\.\d+)?|\.\d+)$'
$cols = @('B', 'C', 'D', 'E')

foreach ($item in $data) {
    $r = $item.Row
    foreach ($col in $cols) {
        if ($item.ContainsKey($col)) {
            $val = $item[$col]
            if ($val -match $numericPattern) {
                $ws.Range("$col$r").Value = "'" + $val
            } else {
                $ws.Range("$col$r").Value = $val
            }
        }
    }
}
